$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$ws.Range("C1").Value = "CaoUser"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C2").Value = "Gemma Hardy"

$ws.Select() | Out-Null
$ws.Range("C1:C2").Select() | Out-Null
